# Append a new record row (row 74) to the "Without Duplication" sheet,
# matching the commit "Update Record 14-04-2025 07:32".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A74").Value = "hm20230396@sva.edu.eg"
$ws.Range("B74").Value = "https://hagersalim.github.io/myfirstweb/"
